$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")
$ws8 = $wb.Worksheets.Item("WVR")

# ALC row 19
$ws1.Range("H19").Value = 5000
$ws1.Range("I19").Value = 5000
$ws1.Range("K19").Value = 5000
$ws1.Range("M19").Value = -4825

# ALC row 62
$ws1.Range("H62").Value = 14854.471
$ws1.Range("I62").Value = 9274.200000000001
$ws1.Range("K62").Value = 9274.200000000001
$ws1.Range("M62").Value = -8650.200000000001

# ALC row 64
$ws1.Range("H64").Value = 69933.336
$ws1.Range("I64").Value = 4898.5
$ws1.Range("K64").Value = 4898.5
$ws1.Range("M64").Value = -4650.5

# ALC row 65
$ws1.Range("H65").Value = 14854.471
$ws1.Range("I65").Value = 9274.200000000001
$ws1.Range("K65").Value = 46371
$ws1.Range("M65").Value = -43251

# ALC row 67
$ws1.Range("H67").Value = 69933.336
$ws1.Range("I67").Value = 4898.5
$ws1.Range("K67").Value = 4898.5
$ws1.Range("M67").Value = -4040.5

# ALC row 106
$ws1.Range("H106").Value = 4367.3335
$ws1.Range("I106").Value = 2796.7896
$ws1.Range("K106").Value = 2796.7896
$ws1.Range("M106").Value = -2165.7896

# ALC row 111
$ws1.Range("H111").Value = 2753.625
$ws1.Range("I111").Value = 2014.5
$ws1.Range("K111").Value = 6043.5
$ws1.Range("M111").Value = -2976.5

# ALC row 113
$ws1.Range("H113").Value = 4336.885
$ws1.Range("I113").Value = 4347.9414
$ws1.Range("J113").Value = 4316
$ws1.Range("K113").Value = 4347.9414
$ws1.Range("L113").Value = 4316
$ws1.Range("M113").Value = -1093.9414
$ws1.Range("N113").Value = -10824

# ALC row 121
$ws1.Range("H121").Value = 1899.5
$ws1.Range("J121").Value = 1899.5
$ws1.Range("L121").Value = 5698.5
$ws1.Range("N121").Value = -9192.5

# ALC row 137
$ws1.Range("H137").Value = 39059.594
$ws1.Range("I137").Value = 57087.39
$ws1.Range("K137").Value = 171262.17
$ws1.Range("M137").Value = -168712.17

# ARM row 32
$ws2.Range("H32").Value = 18127.045
$ws2.Range("I32").Value = 17992.238
$ws2.Range("K32").Value = 17992.238
$ws2.Range("M32").Value = -17705.238

# ARM row 61
$ws2.Range("H61").Value = 13689.579
$ws2.Range("I61").Value = 14039
$ws2.Range("K61").Value = 14039
$ws2.Range("M61").Value = -13827

# ARM row 63
$ws2.Range("H63").Value = 675
$ws2.Range("I63").Value = 675
$ws2.Range("K63").Value = 675
$ws2.Range("M63").Value = 11

# ARM row 66
$ws2.Range("H66").Value = 675
$ws2.Range("I66").Value = 675
$ws2.Range("K66").Value = 3375
$ws2.Range("M66").Value = 57

# ARM row 122
$ws2.Range("H122").Value = 2926.5264
$ws2.Range("I122").Value = 2114.2666
$ws2.Range("K122").Value = 6342.7998
$ws2.Range("M122").Value = -3892.7998

# ARM row 132
$ws2.Range("H132").Value = 25024.445
$ws2.Range("I132").Value = 26356.928
$ws2.Range("J132").Value = 6369.6665
$ws2.Range("K132").Value = 79070.784
$ws2.Range("L132").Value = 19108.9995
$ws2.Range("M132").Value = -76540.784
$ws2.Range("N132").Value = -24168.9995

# ARM row 136
$ws2.Range("H136").Value = 13689.579
$ws2.Range("I136").Value = 14039
$ws2.Range("K136").Value = 42117
$ws2.Range("M136").Value = -39567

# BSM row 105
$ws3.Range("H105").Value = 4026.4546
$ws3.Range("I105").Value = 3258.182
$ws3.Range("J105").Value = 4794.727
$ws3.Range("K105").Value = 3258.182
$ws3.Range("L105").Value = 4794.727
$ws3.Range("M105").Value = -1511.182
$ws3.Range("N105").Value = -8288.726999999999

# BSM row 107
$ws3.Range("H107").Value = 2525.5715
$ws3.Range("I107").Value = 2279.8333
$ws3.Range("J107").Value = 4000
$ws3.Range("K107").Value = 2279.8333
$ws3.Range("L107").Value = 4000
$ws3.Range("M107").Value = -359.8332999999998
$ws3.Range("N107").Value = -7840

# CRP row 22
$ws4.Range("H22").Value = 443.33334
$ws4.Range("I22").Value = 298.2
$ws4.Range("K22").Value = 298.2
$ws4.Range("M22").Value = 51.80000000000001

# CRP row 94
$ws4.Range("H94").Value = 1459.8572
$ws4.Range("J94").Value = 1343.8
$ws4.Range("L94").Value = 1343.8
$ws4.Range("N94").Value = -2245.8

# CRP row 99
$ws4.Range("H99").Value = 22857.4
$ws4.Range("I99").Value = 22857.4
$ws4.Range("K99").Value = 22857.4
$ws4.Range("M99").Value = -21359.4

# CRP row 122
$ws4.Range("H122").Value = 2216.6667
$ws4.Range("I122").Value = 2216.6667
$ws4.Range("J122").Value = 0
$ws4.Range("K122").Value = 6650.000100000001
$ws4.Range("L122").Value = 0
$ws4.Range("N122").Value = -4200.000100000001
$ws4.Range("M122").Value = $null

# CRP row 126
$ws4.Range("H126").Value = 22857.4
$ws4.Range("I126").Value = 22857.4
$ws4.Range("K126").Value = 68572.20000000001
$ws4.Range("M126").Value = -66102.20000000001

# CRP row 132
$ws4.Range("H132").Value = 2753.3333
$ws4.Range("J132").Value = 3299.75
$ws4.Range("L132").Value = 9899.25
$ws4.Range("N132").Value = -14959.25

# CRP row 134
$ws4.Range("H134").Value = 33327.75
$ws4.Range("I134").Value = 37731.785
$ws4.Range("K134").Value = 113195.355
$ws4.Range("M134").Value = -110660.355

# CUL row 5
$ws5.Range("H5").Value = 942
$ws5.Range("I5").Value = 604.75
$ws5.Range("J5").Value = 1391.6666
$ws5.Range("K5").Value = 1814.25
$ws5.Range("L5").Value = 4174.9998
$ws5.Range("M5").Value = -1702.25
$ws5.Range("N5").Value = -4398.9998

# CUL row 12
$ws5.Range("H12").Value = 233.71428
$ws5.Range("I12").Value = 170.25
$ws5.Range("J12").Value = 259.1
$ws5.Range("K12").Value = 510.75
$ws5.Range("L12").Value = 777.3000000000001
$ws5.Range("M12").Value = -337.75
$ws5.Range("N12").Value = -1123.3

# CUL row 38
$ws5.Range("H38").Value = 241.90909
$ws5.Range("I38").Value = 181.71428
$ws5.Range("J38").Value = 347.25
$ws5.Range("K38").Value = 545.14284
$ws5.Range("L38").Value = 1041.75
$ws5.Range("M38").Value = -198.14284
$ws5.Range("N38").Value = -1735.75

# CUL row 131
$ws5.Range("H131").Value = 5560113.5
$ws5.Range("I131").Value = 1600.8182
$ws5.Range("K131").Value = 4802.4546
$ws5.Range("M131").Value = 237.5454

# CUL row 135
$ws5.Range("H135").Value = 942
$ws5.Range("I135").Value = 604.75
$ws5.Range("J135").Value = 1391.6666
$ws5.Range("K135").Value = 5442.75
$ws5.Range("L135").Value = 12524.9994
$ws5.Range("M135").Value = -2907.75
$ws5.Range("N135").Value = -17594.9994

# GSM row 113
$ws6.Range("H113").Value = 136111
$ws6.Range("I113").Value = 93861.45
$ws6.Range("J113").Value = 252297.25
$ws6.Range("K113").Value = 93861.45
$ws6.Range("L113").Value = 252297.25
$ws6.Range("M113").Value = -91691.45
$ws6.Range("N113").Value = -256637.25

# GSM row 132
$ws6.Range("H132").Value = 74648
$ws6.Range("I132").Value = 112841.664
$ws6.Range("K132").Value = 338524.992
$ws6.Range("M132").Value = -335994.992

# GSM row 136
$ws6.Range("H136").Value = 34242.445
$ws6.Range("J136").Value = 34242.445
$ws6.Range("L136").Value = 102727.335
$ws6.Range("N136").Value = -107827.335

# LTW row 7
$ws7.Range("H7").Value = 17120.062
$ws7.Range("I7").Value = 29998.25
$ws7.Range("K7").Value = 29998.25
$ws7.Range("M7").Value = -29886.25

# LTW row 16
$ws7.Range("H16").Value = 3781.1714
$ws7.Range("I16").Value = 4077.8147
$ws7.Range("J16").Value = 2780
$ws7.Range("K16").Value = 4077.8147
$ws7.Range("L16").Value = 2780
$ws7.Range("M16").Value = -3907.8147
$ws7.Range("N16").Value = -3120

# LTW row 22
$ws7.Range("H22").Value = 61677.95
$ws7.Range("I22").Value = 124876.555
$ws7.Range("J22").Value = 4799.2
$ws7.Range("K22").Value = 124876.555
$ws7.Range("L22").Value = 4799.2
$ws7.Range("M22").Value = -124581.555
$ws7.Range("N22").Value = -5389.2

# LTW row 27
$ws7.Range("H27").Value = 61677.95
$ws7.Range("I27").Value = 124876.555
$ws7.Range("J27").Value = 4799.2
$ws7.Range("K27").Value = 124876.555
$ws7.Range("L27").Value = 4799.2
$ws7.Range("M27").Value = -124769.555
$ws7.Range("N27").Value = -5013.2

# LTW row 55
$ws7.Range("H55").Value = 1118.8077
$ws7.Range("I55").Value = 944.1667
$ws7.Range("J55").Value = 1511.75
$ws7.Range("K55").Value = 944.1667
$ws7.Range("L55").Value = 1511.75
$ws7.Range("M55").Value = -771.1667
$ws7.Range("N55").Value = -1857.75

# LTW row 100
$ws7.Range("H100").Value = 3823.8
$ws7.Range("J100").Value = 4061.4285
$ws7.Range("L100").Value = 4061.4285
$ws7.Range("N100").Value = -5143.4285

# LTW row 126
$ws7.Range("H126").Value = 17120.062
$ws7.Range("I126").Value = 29998.25
$ws7.Range("K126").Value = 89994.75
$ws7.Range("M126").Value = -87524.75

# LTW row 132
$ws7.Range("H132").Value = 16970.408
$ws7.Range("I132").Value = 19084.094
$ws7.Range("J132").Value = 5798.0713
$ws7.Range("K132").Value = 57252.28200000001
$ws7.Range("L132").Value = 17394.2139
$ws7.Range("M132").Value = -54722.28200000001
$ws7.Range("N132").Value = -22454.2139

# WVR row 122
$ws8.Range("H122").Value = 128956.5
$ws8.Range("J122").Value = 337066.34
$ws8.Range("L122").Value = 1011199.02
$ws8.Range("N122").Value = -1016099.02

# WVR row 126
$ws8.Range("H126").Value = 85637.60000000001
$ws8.Range("I126").Value = 105671.18
$ws8.Range("K126").Value = 317013.54
$ws8.Range("M126").Value = -314543.54

# WVR row 136
$ws8.Range("H136").Value = 4110.2188
$ws8.Range("I136").Value = 3850.3914
$ws8.Range("K136").Value = 11551.1742
$ws8.Range("M136").Value = -9001.174199999999
